$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header format from G1 (the "sum" header) onto the new H1 header
# cell so the new "Save" column matches the existing header styling.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
